# Update gh-pages to output generated at 456a3b4
# Updates the "views" column (F) for events across sheets:
#   展览 (Exhibitions), 演出 (Performances), 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1027
$ws1.Range("F4").Value  = 5780
$ws1.Range("F5").Value  = 537
$ws1.Range("F6").Value  = 993
$ws1.Range("F7").Value  = 1016
$ws1.Range("F8").Value  = 837
$ws1.Range("F9").Value  = 84
$ws1.Range("F10").Value = 45
$ws1.Range("F11").Value = 605
$ws1.Range("F12").Value = 40
$ws1.Range("F15").Value = 1955
$ws1.Range("F16").Value = 1504
$ws1.Range("F17").Value = 1041
$ws1.Range("F20").Value = 382
$ws1.Range("F21").Value = 605
$ws1.Range("F22").Value = 216
$ws1.Range("F23").Value = 1066
$ws1.Range("F25").Value = 526
$ws1.Range("F26").Value = 3343
$ws1.Range("F32").Value = 456
$ws1.Range("F38").Value = 794
$ws1.Range("F40").Value = 65
$ws1.Range("F42").Value = 79

# ---- Sheet "演出" (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 461
$ws2.Range("F6").Value = 263

# ---- Sheet "全部类型" (All Types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1027
$ws4.Range("F6").Value  = 5780
$ws4.Range("F7").Value  = 537
$ws4.Range("F8").Value  = 993
$ws4.Range("F10").Value = 461
$ws4.Range("F11").Value = 1016
$ws4.Range("F12").Value = 837
$ws4.Range("F14").Value = 263
$ws4.Range("F15").Value = 84
$ws4.Range("F16").Value = 45
$ws4.Range("F17").Value = 605
$ws4.Range("F18").Value = 40
$ws4.Range("F22").Value = 1955
$ws4.Range("F23").Value = 1504
$ws4.Range("F24").Value = 1041
$ws4.Range("F27").Value = 382
$ws4.Range("F29").Value = 605
$ws4.Range("F30").Value = 216
$ws4.Range("F31").Value = 1066
$ws4.Range("F32").Value = 3343
$ws4.Range("F38").Value = 456
$ws4.Range("F43").Value = 794
$ws4.Range("F46").Value = 79
